# "FIN DE 25 SEPT 2021"
# - Bump the "ARQUITECTO" vale amount from $100,000 to $150,000 and update
#   the spelled-out amount text accordingly.
# - Move the active-cell selection on the (active) ARQUITECTO sheet.
# The TODAY() cells on both sheets are volatile formulas and are left
# untouched; they recompute on their own during the workbook's recalc.

$wb = $excel.ActiveWorkbook

$wsArquitecto = $wb.Worksheets.Item("ARQUITECTO        ")

# Numeric amount: 100000 -> 150000
$wsArquitecto.Range("D1").Value = 150000

# Spelled-out amount text that goes with the number above
$wsArquitecto.Range("A2").Value = "CIENTO CINCUENTA     MIL   PESOS 00/100 M.N."

# Re-select the sheet/cell that was active when the file was last saved
[void]$wsArquitecto.Activate()
[void]$wsArquitecto.Range("D3").Select()
